$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prijzen")

# Insert a new row above row 1, shifting all existing data down by one row
$ws.Rows.Item(1).Insert()

# Update the selection to match the new data range
$ws.Range("A2:A8").Select()
